# Work Plan template v0.2
# Builds the weekly task/week grid on Sheet1 and applies the
# "Good" (header) / "20% - Accent1" (body, shaded columns) cell styles,
# matching the look of the Updated Project Work Plan template.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1): Task, Week 1 .. Week 12 -------------------------
$headers = @("Task","Week 1","Week 2","Week 3","Week 4","Week 5","Week 6","Week 7","Week 8","Week 9","Week 10","Week 11","Week 12")
for ($col = 1; $col -le $headers.Length; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# ---- Task names (column A, rows 2-10) -------------------------------------
$tasks = @(
    "Updated Project Work Plan / LR /SRS/ SDO",
    "Test Plan Document",
    "Midterm Interactive-DEMO",
    "Final Product 1st Release",
    "User Manual",
    "Project Report / Test Result / Project Tracking Form",
    "Project Poster",
    "Updated Project Webpage",
    "DEMO Video"
)
for ($i = 0; $i -lt $tasks.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $tasks[$i]
}

# ---- Column A width ---------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 45.140625

# ---- Styling -----------------------------------------------------------
# Header row: built-in "Good" style (green fill)
$ws.Range("A1:M1").Style = "Good"

# Body rows: built-in "20% - Accent1" style (light blue fill) on the
# task-name column plus the shaded "week" columns C, E, G, I, K, M.
foreach ($col in @(1,3,5,7,9,11,13)) {
    $ws.Range($ws.Cells.Item(2, $col), $ws.Cells.Item(10, $col)).Style = "20% - Accent1"
}

# Row 2 also carries the non-shaded columns (B,D,F,H,J,L) explicitly
# (re-applying the default font is a no-op visually, but it registers the
# cell in the sheet so the whole first data row is formatted like the rest).
foreach ($col in @(2,4,6,8,10,12)) {
    $ws.Cells.Item(2, $col).Font.Name = "Calibri"
}

# ---- View / selection ---------------------------------------------------
$ws.Range("P5").Select() | Out-Null

# ---- Page setup ----------------------------------------------------------
$ps = $ws.PageSetup
$ps.Orientation = 2            # xlLandscape
$ps.PaperSize = 9              # xlPaperA4
$ps.LeftMargin = 0.75 * 72
$ps.RightMargin = 0.25 * 72
$ps.TopMargin = 0.75 * 72
$ps.BottomMargin = 0.75 * 72
$ps.HeaderMargin = 0.3 * 72
$ps.FooterMargin = 0.3 * 72
